# "wireframes - ui simple"
# - Fix the risk-description column header typo: "Popis riziko" -> "Popis rizika"
# - Move the sheet's active selection from E10 to B8 (and let the view's
#   top-left/scroll position reset accordingly)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

$ws.Range("B1").Value = "Popis rizika"

$ws.Select()
$ws.Range("B8").Select()
